# Re-applies the latest scraped cryptocurrency market data to the
# "cryptos" worksheet: updated prices / 1h volume-change percentages for
# (almost) every coin row, plus a handful of rows whose coin (name, link)
# pairs were re-ordered because the scrape now ranks them differently.
#
# Column D ("Price") holds values such as "69.826.10" or "0.990" that are
# *text*, not numbers (note the thousands-separator-looking dots and the
# significant trailing zero). Excel's COM layer auto-converts a plain
# numeric-looking string assigned to .Value into a real number, which would
# silently corrupt values like "32.00" -> 32. To keep these as text we
# prefix them with a single quote ( ' ), exactly like typing them by hand
# into a cell in the Excel UI forces text storage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue([string]$cellRef, [string]$text) {
    # Leading apostrophe forces Excel to store the value as text instead of
    # silently parsing it as a number/date.
    $ws.Range($cellRef).Value = "'" + $text
}

Set-TextValue "D2" '69.826.10'
$ws.Range("E2").Value = '  -1.23%  '
Set-TextValue "D3" '3.535.97'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  -0.02%  '
Set-TextValue "D5" '610.63'
$ws.Range("E5").Value = '  +3.56%  '
Set-TextValue "D6" '184.84'
$ws.Range("E6").Value = '  -1.32%  '
$ws.Range("E7").Value = '  -1.40%  '
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +5.06%  '
$ws.Range("E10").Value = '  -1.66%  '
Set-TextValue "D11" '53.57'
$ws.Range("E11").Value = '  -1.88%  '
$ws.Range("E12").Value = '  -1.35%  '
$ws.Range("E13").Value = '  -1.30%  '
Set-TextValue "D14" '4.101.59'
$ws.Range("E14").Value = '  -1.23%  '
$ws.Range("B15").Value = 'BitcoinCash'
$ws.Range("C15").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue "D15" '587.27'
$ws.Range("E15").Value = '  +4.69%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
Set-TextValue "D16" '69.928.67'
$ws.Range("E16").Value = '  -1.07%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue "D17" '3.545.14'
$ws.Range("E17").Value = '  -0.17%  '
$ws.Range("B18").Value = 'Uniswap'
$ws.Range("C18").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
Set-TextValue "D18" '12.64'
$ws.Range("E18").Value = '  +1.20%  '
$ws.Range("B19").Value = 'Chainlink'
$ws.Range("C19").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue "D19" '18.87'
$ws.Range("E19").Value = '  -4.37%  '
$ws.Range("E20").Value = '  -0.40%  '
Set-TextValue "D21" '0.990'
$ws.Range("E21").Value = '  -2.97%  '
$ws.Range("E22").Value = '  -2.14%  '
Set-TextValue "D23" '4.69'
$ws.Range("E23").Value = '  -0.17%  '
Set-TextValue "D24" '4.84'
$ws.Range("E24").Value = '  -1.81%  '
Set-TextValue "D25" '96.38'
$ws.Range("E25").Value = '  +0.60%  '
$ws.Range("E26").Value = '  -0.25%  '
$ws.Range("E27").Value = '  -5.37%  '
Set-TextValue "D28" '9.52'
$ws.Range("E28").Value = '  +3.77%  '
Set-TextValue "D29" '32.00'
$ws.Range("E29").Value = '  -0.89%  '
$ws.Range("E30").Value = '  -4.01%  '
Set-TextValue "D31" '12.14'
$ws.Range("E31").Value = '  -2.88%  '
$ws.Range("E32").Value = '  -1.73%  '
Set-TextValue "D33" '63.29'
$ws.Range("E33").Value = '  -2.82%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue "D34" '3.26'
$ws.Range("E34").Value = '  -1.71%  '
$ws.Range("B35").Value = 'dogwifhat'
$ws.Range("C35").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D35" '3.60'
$ws.Range("E35").Value = '  +17.76%  '
Set-TextValue "D36" '529.27'
$ws.Range("E36").Value = '  -6.31%  '
Set-TextValue "D37" '0.401'
$ws.Range("E37").Value = '  -3.51%  '
$ws.Range("E38").Value = '  -0.03%  '
Set-TextValue "D39" '37.07'
$ws.Range("E39").Value = '  -2.90%  '
$ws.Range("B40").Value = 'Maker'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextValue "D40" '3.529.94'
$ws.Range("E40").Value = '  +5.37%  '
$ws.Range("B41").Value = 'PEPE'
$ws.Range("C41").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
Set-TextValue "D41" '0.0₃0775'
$ws.Range("E41").Value = '  -0.32%  '
$ws.Range("E42").Value = '  +3.60%  '
Set-TextValue "D43" '0.136'
$ws.Range("E43").Value = '  +0.49%  '
$ws.Range("E44").Value = '  +1.98%  '
Set-TextValue "D45" '2.92'
$ws.Range("E45").Value = '  -2.04%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
Set-TextValue "D46" '0.142'
$ws.Range("E46").Value = '  +3.16%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue "D47" '3.36'
$ws.Range("E47").Value = '  -6.21%  '
Set-TextValue "D48" '9.16'
$ws.Range("E48").Value = '  -2.49%  '
$ws.Range("E49").Value = '  +0.32%  '
$ws.Range("E50").Value = '  -4.28%  '
Set-TextValue "D51" '135.47'
$ws.Range("E51").Value = '  -1.27%  '
